# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (interest count) and "最低票价" (lowest price)
# values to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2: 南宁·2024良牙动漫冬季盛典（冬典） — sold out, price column becomes text
    $ws.Range("F2").Value = 10879
    $ws.Range("G2").Value = "已售罄"

    # Row 3: 南宁·第五届小蜜蜂动漫嘉年华
    $ws.Range("F3").Value = 241

    # Row 4: 南宁·0316全职only-全明星周末
    $ws.Range("F4").Value = 71

    # Row 5: 南宁·草莓动漫节
    $ws.Range("F5").Value = 744

    # Row with 南宁·第一届ANE·DACG动漫嘉年华 — row 6 on "展览", row 7 on "全部类型"
    if ($sheetName -eq "展览") {
        $ws.Range("F6").Value = 507
    } else {
        $ws.Range("F7").Value = 507
    }
}
